$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

$ws.Range("C2").Value = 12.8634
$ws.Range("D2").Value = 13.0743
$ws.Range("E2").Value = 27.4406

$ws.Range("C3").Value = 11.1675
$ws.Range("D3").Value = 16.4894
$ws.Range("E3").Value = 23.4498

$ws.Range("C4").Value = 10.767
$ws.Range("D4").Value = 19.5409
$ws.Range("E4").Value = 26.6723

$ws.Range("C5").Value = 9.7902
$ws.Range("D5").Value = 6.8546
$ws.Range("E5").Value = -7.9179

$ws.Range("C6").Value = 9.003299999999999
$ws.Range("D6").Value = 16.6877
$ws.Range("E6").Value = 29.9812

$ws.Range("C7").Value = 8.512
$ws.Range("D7").Value = 10.8059
$ws.Range("E7").Value = 27.1205

$ws.Range("C8").Value = 7.2002
$ws.Range("D8").Value = 12.3843
$ws.Range("E8").Value = 14.5421

$ws.Range("B10").Value = "MARINE"
$ws.Range("C10").Value = 6.3393
$ws.Range("D10").Value = 2.8684
$ws.Range("E10").Value = 15.4303

$ws.Range("B11").Value = "POLICYBZR"
$ws.Range("C11").Value = 6.3363
$ws.Range("D11").Value = 8.714499999999999
$ws.Range("E11").Value = 7.6733

$ws.Range("C12").Value = 6.0394
$ws.Range("D12").Value = 5.0231
$ws.Range("E12").Value = -23.1658

$ws.Range("C13").Value = 5.94
$ws.Range("D13").Value = 9.2964
$ws.Range("E13").Value = 15.4801

$ws.Range("C14").Value = 5.7749
$ws.Range("D14").Value = 6.0594
$ws.Range("E14").Value = 7.7995

$ws.Range("C15").Value = 5.4482
$ws.Range("D15").Value = 11.3233
$ws.Range("E15").Value = 13.4165

$ws.Range("B16").Value = "VENKEYS"
$ws.Range("C16").Value = 5.396
$ws.Range("D16").Value = 6.0253
$ws.Range("E16").Value = 3.7191

$ws.Range("B17").Value = "BLSE"
$ws.Range("C17").Value = 5.2918
$ws.Range("D17").Value = 4.2593
$ws.Range("E17").Value = -1.9269

$ws.Range("C18").Value = 5.2043
$ws.Range("D18").Value = 11.127
$ws.Range("E18").Value = 28.2269

$ws.Range("C23").Value = 4.9726
$ws.Range("D23").Value = 3.5852
$ws.Range("E23").Value = 8.890000000000001

$ws.Range("C24").Value = 4.9513
$ws.Range("D24").Value = 11.5181
$ws.Range("E24").Value = 8.0059

$ws.Range("C25").Value = 4.8578
$ws.Range("D25").Value = 5.0372
$ws.Range("E25").Value = -0.092

$ws.Range("C26").Value = 4.5181
$ws.Range("D26").Value = 4.4132
$ws.Range("E26").Value = 6.0081

$ws.Range("C27").Value = 4.42
$ws.Range("D27").Value = -0.1641
$ws.Range("E27").Value = 38.6515

$ws.Range("B28").Value = "ALICON"
$ws.Range("C28").Value = 4.2891
$ws.Range("D28").Value = 10.4715
$ws.Range("E28").Value = 15.9589

$ws.Range("B29").Value = "MTARTECH"
$ws.Range("C29").Value = 4.2754
$ws.Range("D29").Value = 8.427199999999999
$ws.Range("E29").Value = 32.5267

$ws.Range("C30").Value = 4.1749
$ws.Range("D30").Value = 5.2764
$ws.Range("E30").Value = 6.3677

$ws.Range("C31").Value = 4.1558
$ws.Range("D31").Value = 7.4394
$ws.Range("E31").Value = 0.1613

$ws.Range("B32").Value = "CANBK"
$ws.Range("C32").Value = 4.1314
$ws.Range("D32").Value = 6.6746
$ws.Range("E32").Value = 8.390599999999999

$ws.Range("B33").Value = "BAJAJHCARE"
$ws.Range("C33").Value = 4.0954
$ws.Range("D33").Value = 4.6425
$ws.Range("E33").Value = -1.6365

$ws.Range("B34").Value = "SHANTIGOLD"
$ws.Range("C34").Value = 3.9875
$ws.Range("D34").Value = 11.3065
$ws.Range("E34").Value = 3.8805

$ws.Range("B37").Value = "PFOCUS"
$ws.Range("C37").Value = 3.8493
$ws.Range("D37").Value = 1.1205
$ws.Range("E37").Value = 2.5861

$ws.Range("B38").Value = "SHRINGARMS"
$ws.Range("C38").Value = 3.8469
$ws.Range("D38").Value = 5.0613
$ws.Range("E38").Value = 25.3558

$ws.Range("B39").Value = "BLUEDART"
$ws.Range("C39").Value = 3.7508
$ws.Range("D39").Value = 22.5137
$ws.Range("E39").Value = 19.5599

$ws.Range("B40").Value = "REDTAPE"
$ws.Range("C40").Value = 3.7358
$ws.Range("D40").Value = 3.6351
$ws.Range("E40").Value = -3.3013

$ws.Range("B41").Value = "GMMPFAUDLR"
$ws.Range("C41").Value = 3.6265
$ws.Range("D41").Value = 8.1175
$ws.Range("E41").Value = 20.5341

$ws.Range("B42").Value = "VSTIND"
$ws.Range("C42").Value = 3.597
$ws.Range("D42").Value = 4.1003
$ws.Range("E42").Value = 3.6371

$ws.Range("B43").Value = "SAPPHIRE"
$ws.Range("C43").Value = 3.5855
$ws.Range("D43").Value = 5.412
$ws.Range("E43").Value = 2.757

$ws.Range("B44").Value = "BGRENERGY"
$ws.Range("C44").Value = 3.5346
$ws.Range("D44").Value = -5.8011
$ws.Range("E44").Value = 75.3689

$ws.Range("B45").Value = "RSYSTEMS"
$ws.Range("C45").Value = 3.499
$ws.Range("D45").Value = 4.6608
$ws.Range("E45").Value = 7.011

$ws.Range("B46").Value = "AHLUCONT"
$ws.Range("C46").Value = 3.496
$ws.Range("D46").Value = 2.2983
$ws.Range("E46").Value = -4.9588

$ws.Range("B48").Value = "MIDWESTLTD"
$ws.Range("C48").Value = 3.404
$ws.Range("D48").Value = -1.1837
$ws.Range("E48").Value = "N/A"

$ws.Range("B49").Value = "CENTRUM"
$ws.Range("C49").Value = 3.3636
$ws.Range("D49").Value = 1.9731
$ws.Range("E49").Value = 1.3068

$ws.Range("B50").Value = "RELTD"
$ws.Range("C50").Value = 3.2569
$ws.Range("D50").Value = 10.0472
$ws.Range("E50").Value = -1.4812

$ws.Range("B51").Value = "NEULANDLAB"
$ws.Range("C51").Value = 3.2065
$ws.Range("D51").Value = -1.1686
$ws.Range("E51").Value = 8.8668

$ws.Range("C52").Value = 3.1514
$ws.Range("D52").Value = 3.3261
$ws.Range("E52").Value = 19.033

$ws.Range("B53").Value = "INDORAMA"
$ws.Range("C53").Value = 3.0612
$ws.Range("D53").Value = 5.5381
$ws.Range("E53").Value = 16.9755

$ws.Range("B54").Value = "JKTYRE"
$ws.Range("C54").Value = 3.006
$ws.Range("D54").Value = 6.0566
$ws.Range("E54").Value = 22.1585

$ws.Range("B55").Value = "SUNDROP"
$ws.Range("C55").Value = 2.9786
$ws.Range("D55").Value = 2.8713
$ws.Range("E55").Value = 0.9714

$ws.Range("B56").Value = "PSPPROJECT"
$ws.Range("C56").Value = 2.9047
$ws.Range("D56").Value = 16.9828
$ws.Range("E56").Value = 23.4144

$ws.Range("B57").Value = "FIVESTAR"
$ws.Range("C57").Value = 2.8923
$ws.Range("D57").Value = 15.9311
$ws.Range("E57").Value = 16.0175

$ws.Range("B58").Value = "POWERINDIA"
$ws.Range("C58").Value = 2.8772
$ws.Range("D58").Value = 7.2643
$ws.Range("E58").Value = -0.08890000000000001

$ws.Range("B59").Value = "ASHOKA"
$ws.Range("C59").Value = 2.8461
$ws.Range("D59").Value = 4.3667
$ws.Range("E59").Value = 6.9936

$ws.Range("B60").Value = "INOXGREEN"
$ws.Range("C60").Value = 2.8336
$ws.Range("D60").Value = 10.8511
$ws.Range("E60").Value = 34.2549

$ws.Range("B61").Value = "OIL"
$ws.Range("C61").Value = 2.83
$ws.Range("D61").Value = 3.0751
$ws.Range("E61").Value = 4.4949

$ws.Range("B62").Value = "VOLTAMP"
$ws.Range("C62").Value = 2.8231
$ws.Range("D62").Value = 2.702
$ws.Range("E62").Value = 2.4394

$ws.Range("B63").Value = "BPCL"
$ws.Range("C63").Value = 2.8153
$ws.Range("D63").Value = 8.306800000000001
$ws.Range("E63").Value = 5.3732

$ws.Range("B64").Value = "MRPL"
$ws.Range("C64").Value = 2.7349
$ws.Range("D64").Value = 12.7108
$ws.Range("E64").Value = 23.3376

$ws.Range("B65").Value = "IIFL"
$ws.Range("C65").Value = 2.6769
$ws.Range("D65").Value = 9.6366
$ws.Range("E65").Value = 18.8316

$ws.Range("B66").Value = "CAMS"
$ws.Range("C66").Value = 2.6741
$ws.Range("D66").Value = 2.0205
$ws.Range("E66").Value = 5.3211

$ws.Range("B67").Value = "JKLAKSHMI"
$ws.Range("C67").Value = 2.6028
$ws.Range("D67").Value = 4.5921
$ws.Range("E67").Value = 1.6095

$ws.Range("B69").Value = "MFSL"
$ws.Range("C69").Value = 2.581
$ws.Range("D69").Value = 2.635
$ws.Range("E69").Value = -1.1359

$ws.Range("B70").Value = "REFEX"
$ws.Range("C70").Value = 2.567
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 1.9773

$ws.Range("B71").Value = "SULA"
$ws.Range("C71").Value = 2.5501
$ws.Range("D71").Value = -0.2755
$ws.Range("E71").Value = -2.6999

$ws.Range("C72").Value = 2.5417
$ws.Range("D72").Value = -0.487
$ws.Range("E72").Value = -1.7375

$ws.Range("B73").Value = "FISCHER"
$ws.Range("C73").Value = 2.4707
$ws.Range("D73").Value = 13.4684
$ws.Range("E73").Value = 6.4865

$ws.Range("B74").Value = "KMEW"
$ws.Range("C74").Value = 2.4272
$ws.Range("D74").Value = 4.5365
$ws.Range("E74").Value = -1.0207

$ws.Range("C75").Value = 2.4215
$ws.Range("D75").Value = 5.6585
$ws.Range("E75").Value = 10.2073

$ws.Range("B76").Value = "DBCORP"
$ws.Range("C76").Value = 2.3724
$ws.Range("D76").Value = 5.0019
$ws.Range("E76").Value = 0.9656
